# Actualización automática 2025-09-09 13:25:09
$wb = $excel.ActiveWorkbook

# Sheet: VENTAS POR GRUPO
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M11").Value = 3981.45

# Sheet: VENTA MENSUAL
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F11").Value = 4897.65
$wsMensual.Range("F23").Value = 24202.76

# Sheet: CUMPLIMIENTO MENSUAL
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D12").Value = 21796.16
$wsCumplimiento.Range("E12").Value = 15027.4830921171
$wsCumplimiento.Range("F12").Value = 0.5919066710883352
$wsCumplimiento.Range("D15").Value = 24202.76
$wsCumplimiento.Range("E15").Value = 31221.98316613378
$wsCumplimiento.Range("F15").Value = 0.4366778918118402
